# SRS_Review.xlsx update
# Commit: "Added SRS Version 1.5 including context diagram and modified all
# the requirements in the SRS Updated SRS_review Sheet"
#
# The "Cross review points " sheet has a "Comment" column (I) that records,
# per reviewed requirement, what action was taken in response to the review
# point. Bring those comments up to date with the new wording.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cross review points ")

$ws.Range("I2").Value = "Added status table as required"
$ws.Range("I3").Value = "There was a misconcept but it is modified"
$ws.Range("I4").Value = "Added reference document"
$ws.Range("I5").Value = "Added SW Context diagram descirbing the whole sw signals"
$ws.Range("I6").Value = "modified all the requirements to atomic requirements"
$ws.Range("I7").Value = "modified the requirements to the required table template"
$ws.Range("I8").Value = "modified this requirement to sw not system requirement"
$ws.Range("I9").Value = "modified this requirement to sw not system requirement"

# The "Comment" column wraps text, so the new (longer/shorter) wording
# changes how tall each row needs to be to show it fully.
$ws.Rows.Item(2).RowHeight = 105
$ws.Rows.Item(3).RowHeight = 90
$ws.Rows.Item(4).RowHeight = 90
$ws.Rows.Item(5).RowHeight = 135
$ws.Rows.Item(6).RowHeight = 105
$ws.Rows.Item(7).RowHeight = 120
$ws.Rows.Item(8).RowHeight = 120
$ws.Rows.Item(9).RowHeight = 120
$ws.Rows.Item(10).RowHeight = 30
